$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DOC")

$ws.Range("B2").Value = 5043.814381988261
$ws.Range("C2").Value = 1559.1962603120014
$ws.Range("D2").Value = 12.609535954970655
$ws.Range("E2").Value = 8.406357303313769
$ws.Range("B3").Value = 4842.061806708731
$ws.Range("C3").Value = 1496.8284098995214
$ws.Range("D3").Value = 12.10515451677183
$ws.Range("E3").Value = 8.07010301118122
$ws.Range("B4").Value = 371.1009409278338
$ws.Range("C4").Value = 114.71857516391253
$ws.Range("D4").Value = 0.9277523523195847
$ws.Range("E4").Value = 0.6185015682130565
$ws.Range("B5").Value = 10256.977129624825
$ws.Range("C5").Value = 3170.7432453754354
$ws.Range("D5").Value = 25.64244282406207
$ws.Range("E5").Value = 17.094961882708045
$ws.Range("B7").Value = 4658.228662381298
$ws.Range("D7").Value = 11.645571655953248
$ws.Range("E7").Value = 7.7637144373021645
$ws.Range("B8").Value = 873.4178741964934
$ws.Range("C8").Value = 270.0
$ws.Range("D8").Value = 2.183544685491234
$ws.Range("E8").Value = 1.455696456994156
$ws.Range("B9").Value = 5531.646536577791
$ws.Range("C9").Value = 1710.0
$ws.Range("D9").Value = 13.829116341444482
$ws.Range("E9").Value = 9.21941089429632
$ws.Range("B11").Value = 2024.9926550824753
$ws.Range("C11").Value = 625.9867504718206
$ws.Range("D11").Value = 5.062481637706189
$ws.Range("E11").Value = 3.374987758470793
$ws.Range("B13").Value = 508.2557999999998
$ws.Range("C13").Value = 157.11730896994152
$ws.Range("D13").Value = 1.2706395
$ws.Range("E13").Value = 0.847093
$ws.Range("B14").Value = 422.84368226946447
$ws.Range("C14").Value = 130.71382849564975
$ws.Range("D14").Value = 1.0571092056736615
$ws.Range("E14").Value = 0.7047394704491077
$ws.Range("C15").Value = 409.9068848680968
$ws.Range("D15").Value = 3.315
$ws.Range("E15").Value = 2.21
$ws.Range("C16").Value = 1.6837488291989808
$ws.Range("D16").Value = 0.013616817806294526
$ws.Range("E16").Value = 0.009077878537529683
$ws.Range("C17").Value = 15.367043456626009
$ws.Range("D17").Value = 0.12427639285714281
$ws.Range("E17").Value = 0.08285092857142855
$ws.Range("B18").Value = 2126.7433995348383
$ws.Range("C18").Value = 657.4409968454843
$ws.Range("D18").Value = 5.316858498837098
$ws.Range("E18").Value = 3.5445723325580656
$ws.Range("B20").Value = 1217.9944208403222
$ws.Range("C20").Value = 376.51907905986303
$ws.Range("D20").Value = 3.044986052100806
$ws.Range("E20").Value = 2.0299907014005374
$ws.Range("B21").Value = 4391.925517902091
$ws.Range("C21").Value = 1357.6776075535065
$ws.Range("D21").Value = 10.97981379475523
$ws.Range("E21").Value = 7.319875863170154
$ws.Range("B22").Value = 5720.064789494046
$ws.Range("C22").Value = 1768.2458062633418
$ws.Range("D22").Value = 14.300161973735118
$ws.Range("E22").Value = 9.533441315823412
$ws.Range("B25").Value = 25660.424510313973
$ws.Range("C25").Value = 7932.416798956081
$ws.Range("D25").Value = 64.15106127578495
$ws.Range("E25").Value = 42.76737418385663
$ws.Range("B27").Value = 15403.447380689151
$ws.Range("C27").Value = 4761.673553580647
$ws.Range("D27").Value = 38.50861845172289
$ws.Range("E27").Value = 25.67241230114859